$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 533.55225
$ws.Range("I15").Value = 533.55225
$ws.Range("K15").Value = 1600.65675
$ws.Range("M15").Value = -1431.65675
$ws.Range("H33").Value = 438.16666
$ws.Range("J33").Value = 462.5
$ws.Range("L33").Value = 462.5
$ws.Range("N33").Value = -920.5
$ws.Range("H41").Value = 448.21875
$ws.Range("I41").Value = 247.125
$ws.Range("K41").Value = 247.125
$ws.Range("M41").Value = 192.875
$ws.Range("H51").Value = 14321.75
$ws.Range("I51").Value = 30000
$ws.Range("J51").Value = 9095.666999999999
$ws.Range("K51").Value = 30000
$ws.Range("L51").Value = 9095.666999999999
$ws.Range("M51").Value = -29516
$ws.Range("N51").Value = -10063.667
$ws.Range("H53").Value = 867
$ws.Range("I53").Value = 915.3182
$ws.Range("J53").Value = 791.0714
$ws.Range("K53").Value = 915.3182
$ws.Range("L53").Value = 791.0714
$ws.Range("M53").Value = -278.3182
$ws.Range("N53").Value = -2065.0714
$ws.Range("H76").Value = 101003440
$ws.Range("I76").Value = 168335950
$ws.Range("J76").Value = 4666.25
$ws.Range("K76").Value = 168335950
$ws.Range("L76").Value = 4666.25
$ws.Range("M76").Value = -168335635
$ws.Range("N76").Value = -5296.25
$ws.Range("H79").Value = 101003440
$ws.Range("I79").Value = 168335950
$ws.Range("J79").Value = 4666.25
$ws.Range("K79").Value = 168335950
$ws.Range("L79").Value = 4666.25
$ws.Range("M79").Value = -168334858
$ws.Range("N79").Value = -6850.25
$ws.Range("H98").Value = 6730.6875
$ws.Range("I98").Value = 620.8570999999999
$ws.Range("K98").Value = 620.8570999999999
$ws.Range("M98").Value = 877.1429000000001
$ws.Range("H106").Value = 3371798.2
$ws.Range("I106").Value = 4789187
$ws.Range("K106").Value = 4789187
$ws.Range("M106").Value = -4788556
$ws.Range("H122").Value = 6730.6875
$ws.Range("I122").Value = 620.8570999999999
$ws.Range("K122").Value = 1862.5713
$ws.Range("M122").Value = 587.4287000000002
$ws.Range("H132").Value = 2657.88
$ws.Range("I132").Value = 2415.0454
$ws.Range("J132").Value = 4438.6665
$ws.Range("K132").Value = 7245.1362
$ws.Range("L132").Value = 13315.9995
$ws.Range("M132").Value = -4715.1362
$ws.Range("N132").Value = -18375.9995
$ws.Range("H138").Value = 3209.4082
$ws.Range("I138").Value = 2775.7827
$ws.Range("J138").Value = 3593
$ws.Range("K138").Value = 8327.348100000001
$ws.Range("L138").Value = 10779
$ws.Range("M138").Value = -3187.348100000001
$ws.Range("N138").Value = -21059

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9352.157999999999
$ws.Range("I61").Value = 8109.1816
$ws.Range("K61").Value = 8109.1816
$ws.Range("M61").Value = -7897.1816
$ws.Range("H74").Value = 3630.111
$ws.Range("J74").Value = 6221
$ws.Range("L74").Value = 6221
$ws.Range("N74").Value = -7969
$ws.Range("H77").Value = 3630.111
$ws.Range("J77").Value = 6221
$ws.Range("L77").Value = 31105
$ws.Range("N77").Value = -39841
$ws.Range("H110").Value = 3637.5
$ws.Range("I110").Value = 3592.25
$ws.Range("K110").Value = 3592.25
$ws.Range("M110").Value = -1547.25
$ws.Range("H122").Value = 3837.3333
$ws.Range("I122").Value = 2602
$ws.Range("K122").Value = 7806
$ws.Range("M122").Value = -5356
$ws.Range("H136").Value = 9352.157999999999
$ws.Range("I136").Value = 8109.1816
$ws.Range("K136").Value = 24327.5448
$ws.Range("M136").Value = -21777.5448

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 753.46155
$ws.Range("I64").Value = 457.66666
$ws.Range("J64").Value = 842.2
$ws.Range("K64").Value = 457.66666
$ws.Range("L64").Value = 842.2
$ws.Range("M64").Value = -232.66666
$ws.Range("N64").Value = -1292.2
$ws.Range("H67").Value = 753.46155
$ws.Range("I67").Value = 457.66666
$ws.Range("J67").Value = 842.2
$ws.Range("K67").Value = 457.66666
$ws.Range("L67").Value = 842.2
$ws.Range("M67").Value = 322.33334
$ws.Range("N67").Value = -2402.2
$ws.Range("H107").Value = 1850.7333
$ws.Range("I107").Value = 1841.5385
$ws.Range("J107").Value = 1910.5
$ws.Range("K107").Value = 1841.5385
$ws.Range("L107").Value = 1910.5
$ws.Range("M107").Value = 78.46149999999989
$ws.Range("N107").Value = -5750.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2616.4092
$ws.Range("I31").Value = 1733.9375
$ws.Range("K31").Value = 1733.9375
$ws.Range("M31").Value = -1438.9375
$ws.Range("H34").Value = 2616.4092
$ws.Range("I34").Value = 1733.9375
$ws.Range("K34").Value = 1733.9375
$ws.Range("M34").Value = -1531.9375

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 352.66666
$ws.Range("I8").Value = 352.66666
$ws.Range("K8").Value = 1057.99998
$ws.Range("M8").Value = -918.9999800000001
$ws.Range("H26").Value = 814.125
$ws.Range("J26").Value = 1514.75
$ws.Range("L26").Value = 4544.25
$ws.Range("N26").Value = -5120.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2499
$ws.Range("J122").Value = 2499
$ws.Range("L122").Value = 7497
$ws.Range("N122").Value = -12397

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 113095.22
$ws.Range("I7").Value = 144629.58
$ws.Range("K7").Value = 144629.58
$ws.Range("M7").Value = -144517.58
$ws.Range("H16").Value = 431.14633
$ws.Range("I16").Value = 415.47223
$ws.Range("K16").Value = 415.47223
$ws.Range("M16").Value = -245.47223
$ws.Range("H40").Value = 3329.6667
$ws.Range("I40").Value = 1995
$ws.Range("K40").Value = 1995
$ws.Range("M40").Value = -1859
$ws.Range("H122").Value = 2627.2727
$ws.Range("I122").Value = 2540
$ws.Range("K122").Value = 7620
$ws.Range("M122").Value = -5170
$ws.Range("H126").Value = 113095.22
$ws.Range("I126").Value = 144629.58
$ws.Range("K126").Value = 433888.74
$ws.Range("M126").Value = -431418.74

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 13499.75
$ws.Range("I81").Value = 17333
$ws.Range("K81").Value = 34666
$ws.Range("M81").Value = -33605
$ws.Range("H84").Value = 13499.75
$ws.Range("I84").Value = 17333
$ws.Range("K84").Value = 173330
$ws.Range("M84").Value = -168026
$ws.Range("H107").Value = 1294.2778
$ws.Range("I107").Value = 1294.2778
$ws.Range("K107").Value = 3882.8334
$ws.Range("M107").Value = -1962.8334
$ws.Range("H132").Value = 7726.2915
$ws.Range("I132").Value = 8234.904
$ws.Range("K132").Value = 24704.712
$ws.Range("M132").Value = -22174.712
